# Update "植物大战僵尸2" (pvz2c) download table:
#  - Add a new release row (3.8.8_1725, Nov 6, 2025) at the top of the data (row 2)
#  - Existing rows shift down by one
#  - New row 2 gets hyperlinks on the Baidu Netdisk (D2) and Github (E2) cells

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new blank row above the current first data row (row 2) ---
$ws.Rows("2:2").Insert()

# The row insert shifts everything down, but this engine keeps the old
# hyperlink anchors pinned to the old (now stale) row addresses instead of
# following the shifted cells. Stash the "Aptos Display" hyperlink look
# (now sitting on the shifted C3, formerly C2) before we rebuild hyperlinks,
# so we can restore it after Hyperlinks.Add() resets the font.
$ws.Range("C3").Copy()
$ws.Range("Z1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Clear every stale hyperlink (they still point at the old row-2 addresses)
$ws.Hyperlinks.Delete()

# --- Populate the new row 2 with the newest release info ---
$ws.Range("A2").Value = "3.8.8_1725"
$ws.Range("B2").Value = "Nov 6, 2025"
$ws.Range("C2").Value = "https://pvz2apk-cdn.ditwan.cn/1725/baokai_3.8.8_1725_364_dj2.0-2.0.0.apk"
$ws.Range("D2").Value = "https://pan.baidu.com/s/1-ucUeFy4vhyqC3w795bc1Q?pwd=1234"
$ws.Range("E2").Value = "https://github.com/A-Randomm-User/APK-save/releases/download/pvz2c-tw/3.8.8_1725.apk"

# Row 2 (A2:B2) should look like the other "Version"/"Release Date" cells
$ws.Range("A3:B3").Copy()
$ws.Range("A2:B2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Row 2 (C2) should look like the other plain (non-hyperlinked) Official
# Download Link cells (no hyperlink styling applied there)
$ws.Range("C4").Copy()
$ws.Range("C2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Row 2 should be very slightly taller than the rest
$ws.Rows("2:2").RowHeight = 15

# --- Rebuild hyperlinks in the same order Excel would write them ---
$ws.Hyperlinks.Add($ws.Range("C3"), "https://pvz2apk-cdn.ditwan.cn/1720/baokai_3.8.7_1720_350_dj2.0-2.0.0.apk")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://pan.baidu.com/s/1rM1kpyTDZimR9_fFRtBxmw?pwd=1234")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/A-Randomm-User/APK-save/releases/download/pvz2c-tw/3.8.7_1720.apk")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/A-Randomm-User/APK-save/releases/download/pvz2c-tw/3.8.8_1725.apk")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://pan.baidu.com/s/1-ucUeFy4vhyqC3w795bc1Q?pwd=1234")

# Restore the original "Aptos Display" hyperlink look onto C3:E3 (Hyperlinks.Add
# switches the font to "Aptos Narrow" by default, which is only correct for the
# brand-new row 2 hyperlinks, not the pre-existing row 3 ones).
$ws.Range("Z1").Copy()
$ws.Range("C3:E3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Drop the scratch helper cell entirely so it leaves no trace
$ws.Range("Z1").Clear()

# Match the author's final selection/view state
$ws.Range("C18").Select() | Out-Null
